# Apply updated dSF (column F) values per repulled data / recalculated mean.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = 2
    "F4"  = -1
    "F5"  = -1
    "F6"  = -1
    "F7"  = 4
    "F8"  = -3
    "F10" = 0
    "F12" = 2
    "F14" = -1
    "F15" = 4
    "F17" = 6
    "F18" = 2
    "F19" = 0
    "F20" = 0
    "F22" = 4
    "F29" = -2
    "F30" = -6
    "F31" = -7
    "F32" = 3
    "F33" = -6
    "F34" = 2
    "F35" = -3
    "F36" = 1
    "F37" = 2
    "F39" = -2
    "F41" = -1
    "F44" = -3
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
